$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("N13").ClearContents()

$ws.Range("H15").Value2 = 269.67
$ws.Range("I15").Value2 = 269.67
$ws.Range("K15").Value2 = 809.01
$ws.Range("M15").Value2 = -640.01

$ws.Range("H33").Value2 = 581.4761999999999
$ws.Range("I33").Value2 = 665.85297
$ws.Range("J33").Value2 = 222.875
$ws.Range("K33").Value2 = 665.85297
$ws.Range("L33").Value2 = 222.875
$ws.Range("M33").Value2 = -436.85297
$ws.Range("N33").Value2 = -680.875

$ws.Range("H69").Value2 = 3312.3809
$ws.Range("I69").Value2 = 2483.3333
$ws.Range("J69").Value2 = 3644
$ws.Range("K69").Value2 = 7449.999899999999
$ws.Range("L69").Value2 = 10932
$ws.Range("M69").Value2 = -6575.999899999999
$ws.Range("N69").Value2 = -12680

$ws.Range("H72").Value2 = 3312.3809
$ws.Range("I72").Value2 = 2483.3333
$ws.Range("J72").Value2 = 3644
$ws.Range("K72").Value2 = 22349.9997
$ws.Range("L72").Value2 = 32796
$ws.Range("M72").Value2 = -17981.9997
$ws.Range("N72").Value2 = -41532

$ws.Range("H86").Value2 = 2494.75
$ws.Range("I86").Value2 = 2468.1875
$ws.Range("J86").Value2 = 2601
$ws.Range("K86").Value2 = 2468.1875
$ws.Range("L86").Value2 = 2601
$ws.Range("M86").Value2 = -1345.1875
$ws.Range("N86").Value2 = -4847

$ws.Range("H89").Value2 = 2494.75
$ws.Range("I89").Value2 = 2468.1875
$ws.Range("J89").Value2 = 2601
$ws.Range("K89").Value2 = 12340.9375
$ws.Range("L89").Value2 = 13005
$ws.Range("M89").Value2 = -6724.9375
$ws.Range("N89").Value2 = -24237

$ws.Range("H106").Value2 = 1994.1666
$ws.Range("I106").Value2 = 1994.1666
$ws.Range("K106").Value2 = 1994.1666
$ws.Range("M106").Value2 = -1363.1666

$ws.Range("H132").Value2 = 2912.2307
$ws.Range("I132").Value2 = 2634.116
$ws.Range("J132").Value2 = 5044.4443
$ws.Range("K132").Value2 = 7902.348
$ws.Range("L132").Value2 = 15133.3329
$ws.Range("M132").Value2 = -5372.348
$ws.Range("N132").Value2 = -20193.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1757.5094
$ws.Range("I61").Value2 = 1504.5349
$ws.Range("K61").Value2 = 1504.5349
$ws.Range("M61").Value2 = -1292.5349

$ws.Range("H136").Value2 = 1757.5094
$ws.Range("I136").Value2 = 1504.5349
$ws.Range("K136").Value2 = 4513.6047
$ws.Range("M136").Value2 = -1963.6047

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value2 = 25847.5
$ws.Range("J50").Value2 = 25847.5
$ws.Range("L50").Value2 = 25847.5
$ws.Range("N50").Value2 = -26995.5

$ws.Range("H94").Value2 = 13208.92
$ws.Range("I94").Value2 = 10111.909
$ws.Range("K94").Value2 = 10111.909
$ws.Range("M94").Value2 = -9660.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 2110.628
$ws.Range("I31").Value2 = 1349.1
$ws.Range("J31").Value2 = 2772.8262
$ws.Range("K31").Value2 = 1349.1
$ws.Range("L31").Value2 = 2772.8262
$ws.Range("M31").Value2 = -1054.1
$ws.Range("N31").Value2 = -3362.8262

$ws.Range("H34").Value2 = 2110.628
$ws.Range("I34").Value2 = 1349.1
$ws.Range("J34").Value2 = 2772.8262
$ws.Range("K34").Value2 = 1349.1
$ws.Range("L34").Value2 = 2772.8262
$ws.Range("M34").Value2 = -1147.1
$ws.Range("N34").Value2 = -3176.8262

$ws.Range("H122").Value2 = 1211.9524
$ws.Range("I122").Value2 = 950.0769
$ws.Range("J122").Value2 = 1637.5
$ws.Range("K122").Value2 = 2850.2307
$ws.Range("L122").Value2 = 4912.5
$ws.Range("M122").Value2 = -400.2307000000001
$ws.Range("N122").Value2 = -9812.5

$ws.Range("H132").Value2 = 1897.7273
$ws.Range("I132").Value2 = 1707.8948
$ws.Range("J132").Value2 = 3100
$ws.Range("K132").Value2 = 5123.6844
$ws.Range("L132").Value2 = 9300
$ws.Range("M132").Value2 = -2593.6844
$ws.Range("N132").Value2 = -14360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 809.1905
$ws.Range("J5").Value2 = 1213
$ws.Range("L5").Value2 = 3639
$ws.Range("N5").Value2 = -3863

$ws.Range("H68").Value2 = 1505
$ws.Range("I68").Value2 = 425
$ws.Range("J68").Value2 = 2045
$ws.Range("K68").Value2 = 1275
$ws.Range("L68").Value2 = 6135
$ws.Range("M68").Value2 = -464
$ws.Range("N68").Value2 = -7757

$ws.Range("H71").Value2 = 1505
$ws.Range("I71").Value2 = 425
$ws.Range("J71").Value2 = 2045
$ws.Range("K71").Value2 = 3825
$ws.Range("L71").Value2 = 18405
$ws.Range("M71").Value2 = 231
$ws.Range("N71").Value2 = -26517

$ws.Range("H113").Value2 = 1104.2
$ws.Range("I113").Value2 = 472.4375
$ws.Range("J113").Value2 = 1983.174
$ws.Range("K113").Value2 = 1417.3125
$ws.Range("L113").Value2 = 5949.522
$ws.Range("M113").Value2 = 752.6875
$ws.Range("N113").Value2 = -10289.522

$ws.Range("H132").Value2 = 1158.1786
$ws.Range("I132").Value2 = 837
$ws.Range("J132").Value2 = 1366
$ws.Range("K132").Value2 = 7533
$ws.Range("L132").Value2 = 12294
$ws.Range("M132").Value2 = -5003
$ws.Range("N132").Value2 = -17354

$ws.Range("H135").Value2 = 809.1905
$ws.Range("J135").Value2 = 1213
$ws.Range("L135").Value2 = 10917
$ws.Range("N135").Value2 = -15987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value2 = 40031.5
$ws.Range("I134").Value2 = 40000
$ws.Range("J134").Value2 = 40042
$ws.Range("K134").Value2 = 120000
$ws.Range("L134").Value2 = 120126
$ws.Range("M134").Value2 = -117465
$ws.Range("N134").Value2 = -125196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 2083980.5
$ws.Range("I22").Value2 = 4762124.5
$ws.Range("J22").Value2 = 979.6667
$ws.Range("K22").Value2 = 4762124.5
$ws.Range("L22").Value2 = 979.6667
$ws.Range("M22").Value2 = -4761829.5
$ws.Range("N22").Value2 = -1569.6667

$ws.Range("H27").Value2 = 2083980.5
$ws.Range("I27").Value2 = 4762124.5
$ws.Range("J27").Value2 = 979.6667
$ws.Range("K27").Value2 = 4762124.5
$ws.Range("L27").Value2 = 979.6667
$ws.Range("M27").Value2 = -4762017.5
$ws.Range("N27").Value2 = -1193.6667

$ws.Range("H55").Value2 = 444.66666
$ws.Range("I55").Value2 = 518
$ws.Range("J55").Value2 = 371.33334
$ws.Range("K55").Value2 = 518
$ws.Range("L55").Value2 = 371.33334
$ws.Range("M55").Value2 = -345
$ws.Range("N55").Value2 = -717.33334

$ws.Range("H93").Value2 = 2151.0908
$ws.Range("I93").Value2 = 2815.6
$ws.Range("J93").Value2 = 1597.3334
$ws.Range("K93").Value2 = 2815.6
$ws.Range("L93").Value2 = 1597.3334
$ws.Range("M93").Value2 = -1567.6
$ws.Range("N93").Value2 = -4093.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value2 = 0
$ws.Range("J12").Value2 = 0
$ws.Range("L12").Value2 = 0
$ws.Range("N12").ClearContents()

$ws.Range("H62").Value2 = 171227.78
$ws.Range("I62").Value2 = 4862.5
$ws.Range("J62").Value2 = 304320
$ws.Range("K62").Value2 = 4862.5
$ws.Range("L62").Value2 = 304320
$ws.Range("M62").Value2 = -4238.5
$ws.Range("N62").Value2 = -305568

$ws.Range("H65").Value2 = 171227.78
$ws.Range("I65").Value2 = 4862.5
$ws.Range("J65").Value2 = 304320
$ws.Range("K65").Value2 = 24312.5
$ws.Range("L65").Value2 = 1521600
$ws.Range("M65").Value2 = -21192.5
$ws.Range("N65").Value2 = -1527840

$ws.Range("H81").Value2 = 83336330
$ws.Range("I81").Value2 = 333337300
$ws.Range("J81").Value2 = 2672.2222
$ws.Range("K81").Value2 = 666674600
$ws.Range("L81").Value2 = 5344.4444
$ws.Range("M81").Value2 = -666673539
$ws.Range("N81").Value2 = -7466.4444

$ws.Range("H84").Value2 = 83336330
$ws.Range("I84").Value2 = 333337300
$ws.Range("J84").Value2 = 2672.2222
$ws.Range("K84").Value2 = 3333373000
$ws.Range("L84").Value2 = 26722.222
$ws.Range("M84").Value2 = -3333367696
$ws.Range("N84").Value2 = -37330.222

$ws.Range("H126").Value2 = 3058.077
$ws.Range("I126").Value2 = 5846.3335
$ws.Range("J126").Value2 = 668.1429000000001
$ws.Range("K126").Value2 = 17539.0005
$ws.Range("L126").Value2 = 2004.4287
$ws.Range("M126").Value2 = -15069.0005
$ws.Range("N126").Value2 = -6944.4287
